{"js": "// Apply the two substantive content edits captured by the diff:\n//  1) Remove the parenthetical \"(section \u201cTaxa\u201d)\" from the \"Taxonomy\" bullet,\n//     leaving a Word-style \"_GoBack\" bookmark at the point of the last edit.\n//  2) Replace \"Renku\" with a quoted path \"Documents/Bio334_Data\" (set in\n//     Menlo 11pt black, matching a pasted-from-Terminal style) in the final\n//     bullet about uploading data.\n\n// --- Edit 1: drop \"(section \u201cTaxa\u201d)\" ------------------------------------\nconst taxaResults = context.document.body.search(\n  \"(section \\u201cTaxa\\u201d),\",\n  { matchCase: true }\n);\ntaxaResults.load(\"text\");\nawait context.sync();\n\nif (taxaResults.items.length > 0) {\n  const taxaRange = taxaResults.items[0];\n  // Keep just the comma; the rest of the parenthetical goes away.\n  const commaRange = taxaRange.insertText(\",\", Word.InsertLocation.replace);\n  await context.sync();\n\n  // Word stamps the last edit location with a collapsed \"_GoBack\" bookmark.\n  const collapsedEnd = commaRange.getRange(Word.RangeLocation.end);\n  await context.sync();\n  collapsedEnd.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// --- Edit 2: \"Renku\" -> quoted \"Documents/Bio334_Data\" path -------------\nconst renkuResults = context.document.body.search(\"Renku\", { matchCase: true });\nrenkuResults.load(\"text\");\nawait context.sync();\n\nif (renkuResults.items.length > 0) {\n  const renkuRange = renkuResults.items[0];\n  // Insert the whole replacement as one run first so it inherits the\n  // surrounding (Arial) formatting for the quote characters.\n  renkuRange.insertText('\"Documents/Bio334_Data\" ', Word.InsertLocation.replace);\n  await context.sync();\n\n  // Now re-find just the path portion and give it the monospace styling.\n  const pathResults = context.document.body.search(\"Documents/Bio334_Data\", {\n    matchCase: true,\n  });\n  pathResults.load(\"text\");\n  await context.sync();\n\n  if (pathResults.items.length > 0) {\n    const pathRange = pathResults.items[0];\n    pathRange.font.name = \"Menlo\";\n    pathRange.font.size = 11;\n    pathRange.font.color = \"#000000\";\n    await context.sync();\n  }\n}\n", "ps1": "# Apply the two substantive content edits captured by the diff:\n#  1) Remove the parenthetical \"(section \u201cTaxa\u201d)\" from the \"Taxonomy\" bullet,\n#     leaving a Word-style \"_GoBack\" bookmark at the point of the last edit.\n#  2) Replace \"Renku\" with a quoted path \"Documents/Bio334_Data\" (set in\n#     Menlo 11pt black, matching a pasted-from-Terminal style) in the final\n#     bullet about uploading data.\n\n$doc = $word.ActiveDocument\n\n# --- Edit 1: drop \"(section \u201cTaxa\u201d)\" --------------------------------------\n$openQuote = [char]8220\n$closeQuote = [char]8221\n\n$taxaRange = $doc.Content\n$taxaFind = $taxaRange.Find\n$found = $taxaFind.Execute(\"(section \" + $openQuote + \"Taxa\" + $closeQuote + \"),\")\n\nif ($found) {\n    # Keep just the comma; the rest of the parenthetical goes away.\n    $taxaRange.Text = \",\"\n    $taxaRange.Collapse(0)\n    $doc.Bookmarks.Add(\"_GoBack\", $taxaRange)\n}\n\n# --- Edit 2: \"Renku\" -> quoted \"Documents/Bio334_Data\" path --------------\n$renkuRange = $doc.Content\n$renkuFind = $renkuRange.Find\n$foundRenku = $renkuFind.Execute(\"Renku\")\n\nif ($foundRenku) {\n    # Insert the whole replacement as one run first so it inherits the\n    # surrounding (Arial) formatting for the quote characters.\n    $renkuRange.Text = [char]34 + \"Documents/Bio334_Data\" + [char]34 + \" \"\n\n    # Now re-find just the path portion and give it the monospace styling.\n    $pathRange = $doc.Content\n    $pathFind = $pathRange.Find\n    $foundPath = $pathFind.Execute(\"Documents/Bio334_Data\")\n\n    if ($foundPath) {\n        $pathRange.Font.Name = \"Menlo\"\n        $pathRange.Font.Size = 11\n        $pathRange.Font.Color = 0\n    }\n}\n"}
